# multiple_session_batch.xlsx - "batch definition updates again"
#
# Changes to the "Sessions" sheet:
#  1. Insert two new rows right after "Batch Name" (row 3):
#       - "Context Folder Name" (String)
#       - "Context Name" (String)
#     (row 3's blank spacer + "Session Settings" header shift down to make
#      room, ending up two rows lower than before)
#  2. Remove the "Session Output Folder Name" row from the Session Settings
#     block.
#  3. Remove the "Database Dump Folder Name" row from the IO Settings block.
#  4. Update the frozen-pane / selection to match the new layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Insert two blank rows before the old row 4 -----------------------
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).Insert()

# Populate the two new rows (Context Folder Name / Context Name), mirroring
# the "Batch Name" row (row 3) layout: A = label, B = "String", C = blank.
$ws.Range("A4").Value = "Context Folder Name"
$ws.Range("B4").Value = "String"

$ws.Range("A5").Value = "Context Name"
$ws.Range("B5").Value = "String"

# Copy the formatting from row 3 (A:C) down onto the two new rows so the
# cell styles match the rest of the "label / type / value" rows.
$ws.Range("A3:C3").Copy()
$ws.Range("A4:C4").PasteSpecial(-4122)
$ws.Range("A3:C3").Copy()
$ws.Range("A5:C5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 2 & 3. Remove "Database Dump Folder Name" and "Session Output Folder
#            Name" rows. After the inserts above, the old row 8 ("Session
#            Output Folder Name") now sits at row 10, and the old row 11
#            ("Database Dump Folder Name") now sits at row 13. Delete the
#            higher-numbered row first so the lower index stays valid.
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(10).Delete()

# --- 4. Update the frozen pane / selection --------------------------------
# Re-freeze two rows lower (ySplit 7->9, top-left cell A8->A10) now that the
# two new rows pushed the "Session Settings" block down.
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A10").Select()
$win.FreezePanes = $true
$ws.Range("C5").Select()
